# Update the Chai Tea market-trends header labels (table header row / shared
# strings) to the new wording from the commit. All four changed headers keep
# their existing bold + white Calibri look - B1/C1/D1 were already bold+white
# and E1 becomes bold+white to match the others - so we restore that font
# after updating each cell's text. Updating the cell values also keeps the
# backing Excel Table's column names (Table1) in sync automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1: "チャイの売上合計 (単位)" -> "チャイの売上合計 (ユニット数)"
$b1 = $ws.Range("B1")
$b1.Value = "チャイの売上合計 (ユニット数)"
$b1.Font.Bold = $true
$b1.Font.Color = 16777215

# C1: "職人チャイ販売 (ユニット)" -> "Artisanal Chai の販売 (ユニット数)"
$c1 = $ws.Range("C1")
$c1.Value = "Artisanal Chai の販売 (ユニット数)"
$c1.Font.Bold = $true
$c1.Font.Color = 16777215

# D1: "事前に作成されたチャイの売上 (単位)" -> "事前に作成されたチャイの売上 (ユニット数)"
$d1 = $ws.Range("D1")
$d1.Value = "事前に作成されたチャイの売上 (ユニット数)"
$d1.Font.Bold = $true
$d1.Font.Color = 16777215

# E1: "ソーシャル メディア エンゲージメント" -> "ソーシャル メディア エンゲージメント (ビュー)"
# also gains bold (it previously lacked it while the other headers had it).
$e1 = $ws.Range("E1")
$e1.Value = "ソーシャル メディア エンゲージメント (ビュー)"
$e1.Font.Bold = $true
$e1.Font.Color = 16777215
